# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet
#    that carries it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#    Overview columns E:F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text wherever it appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Shrink the Status column widths ---
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
